$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous contents (old layout was A1:D3)
$ws.Range("A1:D3").Clear()

# New header row: stimulus array (stim1) separated from probe
$ws.Range("A1").Value = "stim1_x"
$ws.Range("B1").Value = "stim1_y"
$ws.Range("C1").Value = "sitm1_color"
$ws.Range("D1").Value = "probe1_x"
$ws.Range("E1").Value = "probe1_y"
$ws.Range("F1").Value = "probe1_color"
$ws.Range("G1").Value = "answer"

# New data row: fixations added for stim1 + probe1, plus the answer/probe1 outcome
$ws.Range("A2").Value = -0.25
$ws.Range("B2").Value = -0.25
$ws.Range("C2").Value = "black"
$ws.Range("D2").Value = -0.25
$ws.Range("E2").Value = -0.25
$ws.Range("F2").Value = "white"
$ws.Range("G2").Value = "probe1"

$ws.Columns.Item(6).ColumnWidth = 11

$ws.Range("B11").Select()
